# Updated symbol list on Mon Dec 12 08:46:53 UTC 2022 with GitHub Actions
#
# This script applies the latest price/volume refresh to the crypto
# tracker sheet. The "Price" (column D) values are stored as text, so we
# force a text number format before writing them back to avoid Excel
# auto-converting the numeric-looking strings into floating point
# numbers (which would lose trailing/leading zero formatting).
#
# Additionally, two rows swapped rank position: "KickToken" (was row 41)
# and "CEJI" (was row 42) traded places, each picking up a freshly
# refreshed price/volume value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$Text
    )
    $range = $ws.Range($CellAddress)
    $range.NumberFormat = "@"
    $range.Value = $Text
}

# --- Price (column D) refreshes -----------------------------------------
Set-TextValue "D2"  "281.38"
Set-TextValue "D3"  "20.81"
Set-TextValue "D4"  "6.240"
Set-TextValue "D5"  "0.06141"
Set-TextValue "D6"  "3.574"
Set-TextValue "D7"  "6.562"
Set-TextValue "D8"  "1.467"
Set-TextValue "D9"  "0.8171"
Set-TextValue "D11" "0.1629"
Set-TextValue "D12" "0.08313"
Set-TextValue "D13" "0.03531"
Set-TextValue "D14" "0.03207"
Set-TextValue "D15" "0.09148"
Set-TextValue "D16" "3.724"
Set-TextValue "D17" "0.001639"
Set-TextValue "D18" "0.04653"
Set-TextValue "D19" "0.006437"
Set-TextValue "D20" "0.006164"
Set-TextValue "D21" "0.001067"
Set-TextValue "D23" "3.806"
Set-TextValue "D24" "2.337"
Set-TextValue "D25" "0.3335"
Set-TextValue "D40" "0.04662"

# --- Rows 41/42 swap: KickToken <-> CEJI ---------------------------------
Set-TextValue "B41" "CEJI"
Set-TextValue "C41" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D41" "0.003512"
Set-TextValue "E41" "40CEJICEJI"

Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.007171"
Set-TextValue "E42" "41KickTokenKICK"

# --- Remaining price refreshes -------------------------------------------
Set-TextValue "D43" "0.1100"
Set-TextValue "D44" "0.01133"
Set-TextValue "D45" "0.00006367"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D48" "0.002944"
Set-TextValue "D49" "0.00001902"
